# Weekly fruit/hortaliza update: insert a new weekly record row at row 55
# (shifting all subsequent rows down by one), matching the data that
# row 55 previously held except for the Fecha (date) and Volumen values,
# which carry the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55; this shifts rows 55..132 down to 56..133
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly entry.
$ws.Cells.Item(55, 1).Value2 = 4
$ws.Cells.Item(55, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(55, 3).Value2 = "Los Lagos"
$ws.Cells.Item(55, 4).Value2 = 44799
$ws.Cells.Item(55, 5).Value2 = 10
$ws.Cells.Item(55, 6).Value2 = 100112022
$ws.Cells.Item(55, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(55, 8).Value2 = "Perfection"
$ws.Cells.Item(55, 9).Value2 = "Primera"
$ws.Cells.Item(55, 10).Value2 = 70
$ws.Cells.Item(55, 11).Value2 = 46000
$ws.Cells.Item(55, 12).Value2 = 46000
$ws.Cells.Item(55, 13).Value2 = 46000
$ws.Cells.Item(55, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(55, 15).Value2 = "Provincia de Huasco"
$ws.Cells.Item(55, 16).Value2 = 1840
$ws.Cells.Item(55, 17).Value2 = 25
$ws.Cells.Item(55, 18).Value2 = "Hortaliza"
